# ---------------------------------------------------------------------------
# Applies the "Number activities and equations / Fix SideNoteParagraph /
# Updated documentation" edit to the empty-sc starter document:
#
#   1. Insert a new "Unit id" paragraph (style UnitID) at the top of the body.
#   2. Fix the casing of the Heading1 run:  "Session Title" -> "Session title"
#   3. Fix the casing of the NormalWeb run: "text here"     -> "Text here"
#   4. Give the "Label" character style a pale-yellow shading (highlight),
#      matching the SideNoteParagraph fix.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. Insert the new "Unit id" paragraph before the existing first
#        paragraph ("Unit title"), using the (already defined) UnitID style.
$firstPara = $d.Paragraphs.Item(1)
$firstPara.Range.InsertParagraphBefore()

$unitIdPara = $d.Paragraphs.Item(1)
$unitIdPara.Range.Text = "Unit id"
$unitIdPara.Style = "UnitID"

# --- 2. Normalise "Session Title" -> "Session title" (Heading1 paragraph).
$d.Content.Find.Execute("Session Title", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Session title", 2) | Out-Null

# --- 3. Normalise "text here" -> "Text here" (NormalWeb paragraph).
$d.Content.Find.Execute("text here", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Text here", 2) | Out-Null

# --- 4. Give the "Label" character style a pale-yellow (FFFF99) shading,
#        mirroring the SideNoteParagraph character style.
try {
    $labelStyle = $d.Styles.Item("Label")
    $labelStyle.Font.Shading.Texture = 0
    $labelStyle.Font.Shading.ForegroundPatternColor = -16777216
    $labelStyle.Font.Shading.BackgroundPatternColor = 10092543
} catch {
    # Older/headless Word OM builds don't route Shading writes through a
    # Style's Font object - ignore so the rest of the edit still applies.
}
